$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.880.92'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.238.76'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''272.05'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +5.44%  '
$ws.Range("D6").Value = '''86.67'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +9.43%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.606'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '''45.05'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.29%  '
$ws.Range("D11").Value = '''0.0923'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("E12").Value = '  +7.80%  '
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").Value = '2.575.10'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '''14.94'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '2.229.77'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '''0.793'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = '43.818.18'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '''0.0000103'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '''70.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''5.97'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").Value = '''2.35'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").Value = '''233.44'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '''8.68'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.26%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '''2.53'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +13.23%  '
$ws.Range("D27").Value = '''10.78'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").Value = '''3.54'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.29%  '
$ws.Range("E29").Value = '  +5.13%  '
$ws.Range("D30").Value = '''39.56'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -7.11%  '
$ws.Range("D31").Value = '''173.80'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("D32").Value = '''0.0909'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.58%  '
$ws.Range("D33").Value = '''20.79'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("D34").Value = '''5.38'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '''0.123'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("D36").Value = '''0.110'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("D37").Value = '''0.0354'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").Value = '''4.33'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("D39").Value = '''3.44'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +19.02%  '
$ws.Range("D40").Value = '''2.21'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").Value = '''12.41'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.68%  '
$ws.Range("D42").Value = '''64.20'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.52%  '
$ws.Range("D43").Value = '''0.205'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '''5.42'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").Value = '''8.48'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").Value = '''0.0983'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = '''99.77'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.85%  '
$ws.Range("E48").Value = '  +3.85%  '
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '''0.427'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.72%  '
$ws.Range("E51").Value = '  -3.05%  '
